$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.972.52"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.556.93"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.33"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.17"
$ws.Range("E8").Value = "  +4.29%  "
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0858"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.778.75"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.556.80"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("E14").Value = "  +1.47%  "
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.973.02"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.83"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "218.70"
$ws.Range("E18").Value = "  +2.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0698"
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.33"
$ws.Range("E20").Value = "  +2.13%  "
$ws.Range("E22").Value = "  +1.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.26"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.25"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.65"
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("E27").Value = "  +1.32%  "
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  +2.48%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.426.62"
$ws.Range("E33").Value = "  +5.24%  "
$ws.Range("E34").Value = "  +5.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  +3.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.981"
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.522"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.75"
$ws.Range("E41").Value = "  +3.00%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("E43").Value = "  +5.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.988"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.57"
$ws.Range("E45").Value = "  +1.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.76"
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.692.48"
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.01"
$ws.Range("E48").Value = "  +2.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0520"
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₇0999"
$ws.Range("E50").Value = "  +3.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0958"
$ws.Range("E51").Value = "  +1.17%  "
